$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 2 (A2:T2) -- Target cluster becomes "ECs", values recomputed
# -----------------------------------------------------------------
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.676920666666666
$ws.Range("N2").Value = 8.030761999999999
$ws.Range("O2").Value = 0.1137411923116975
$ws.Range("P2").Value = 0.1167596329733683
$ws.Range("Q2").Value = 1.039719556160889
$ws.Range("R2").Value = 9.357476005448
$ws.Range("S2").Value = 0.1137411923116975
$ws.Range("T2").Value = 0.1167596329733683

# -----------------------------------------------------------------
# Row 3 (A3:T3) -- Target cluster becomes "FAPs", values recomputed
# -----------------------------------------------------------------
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.75807
$ws.Range("N3").Value = 32.27421
$ws.Range("O3").Value = 0.4571057050773153
$ws.Range("P3").Value = 0.4692362834442626
$ws.Range("Q3").Value = 4.178448732093333
$ws.Range("R3").Value = 37.60603858884
$ws.Range("S3").Value = 0.4571057050773153
$ws.Range("T3").Value = 0.4692362834442626

# -----------------------------------------------------------------
# Row 4 (A4:T4) -- Target cluster becomes "M1", values recomputed
# -----------------------------------------------------------------
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.564735333333333
$ws.Range("N4").Value = 10.694206
$ws.Range("O4").Value = 0.1514640505181089
$ws.Range("P4").Value = 0.1554835727296604
$ws.Range("Q4").Value = 1.384547956447111
$ws.Range("R4").Value = 12.460931608024
$ws.Range("S4").Value = 0.1514640505181089
$ws.Range("T4").Value = 0.1554835727296604

# -----------------------------------------------------------------
# Row 5 (A5:T5) -- Target cluster becomes "M2", values recomputed
# -----------------------------------------------------------------
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.710187333333333
$ws.Range("N5").Value = 14.130562
$ws.Range("O5").Value = 0.2001338067190093
$ws.Range("P5").Value = 0.2054449170361947
$ws.Range("Q5").Value = 1.829443040516444
$ws.Range("R5").Value = 16.464987364648
$ws.Range("S5").Value = 0.2001338067190093
$ws.Range("T5").Value = 0.2054449170361947

# -----------------------------------------------------------------
# Row 6 (A6:T6) -- NEW row, Target cluster "sCs"
# -----------------------------------------------------------------
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd5"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3884013333333334
$ws.Range("H6").Value = 1.165204
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.8252775
$ws.Range("N6").Value = 3.650555
$ws.Range("O6").Value = 0.07755524537386906
$ws.Range("P6").Value = 0.05307559381651387
$ws.Range("Q6").Value = 0.7089402147033333
$ws.Range("R6").Value = 4.25364128822
$ws.Range("S6").Value = 0.07755524537386906
$ws.Range("T6").Value = 0.05307559381651387
